$wb = $excel.ActiveWorkbook

# Update OFF sheet (row 3: Road totals) with Week 16 logged + season sim from Week 17
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 327
$wsOff.Range("C3").Value = 213
$wsOff.Range("D3").Value = 145
$wsOff.Range("E3").Value = 57

# Update DEF sheet (row 3: Road totals) with Week 16 logged + season sim from Week 17
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 444
$wsDef.Range("C3").Value = 338
$wsDef.Range("D3").Value = 106
$wsDef.Range("E3").Value = 52
